# Atualizado por script em 20-12-2023 14:45
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns F..V hold the match data (home/away teams, scores, odds, timestamps,
# url) for each fixture row; A..E (index/country/tournament/season/date) stay
# put. A handful of fixture rows had their match-data swapped with their
# neighbour row - fix them by swapping columns F:V between the row pairs.
function Swap-MatchData($row1, $row2) {
    $cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
    foreach ($col in $cols) {
        $c1 = $ws.Range("$col$row1")
        $c2 = $ws.Range("$col$row2")
        $v1 = $c1.Value()
        $v2 = $c2.Value()
        $c1.Value = $v2
        $c2.Value = $v1
    }
}

Swap-MatchData 63 64
Swap-MatchData 68 69
Swap-MatchData 87 88
Swap-MatchData 89 90
Swap-MatchData 112 113

# Append the new fixture as row 154, reusing the formatting of row 153
# (bold/bordered/centered index column, datetime-formatted match-date column).
$ws.Range("A153").Copy()
$ws.Range("A154").PasteSpecial(-4122)
$ws.Range("E153").Copy()
$ws.Range("E154").PasteSpecial(-4122)

$ws.Range("A154").Value = 153
$ws.Range("B154").Value = "turkey"
$ws.Range("C154").Value = "super-lig"
$ws.Range("D154").Value = "2023-2024"
$ws.Range("E154").Value = 45280.625
$ws.Range("F154").Value = "Kayserispor"
$ws.Range("G154").Value = 3
$ws.Range("H154").Value = "Fenerbahce"
$ws.Range("I154").Value = 4
$ws.Range("J154").Value = 7.02
$ws.Range("K154").Value = "14/12/2023 09:48"
$ws.Range("L154").Value = 5.73
$ws.Range("M154").Value = "20/12/2023 14:59"
$ws.Range("N154").Value = 5.17
$ws.Range("O154").Value = "14/12/2023 09:48"
$ws.Range("P154").Value = 4.48
$ws.Range("Q154").Value = "20/12/2023 14:59"
$ws.Range("R154").Value = 1.41
$ws.Range("S154").Value = "14/12/2023 09:48"
$ws.Range("T154").Value = 1.58
$ws.Range("U154").Value = "20/12/2023 14:59"
$ws.Range("V154").Value = "https://www.betexplorer.com/football/turkey/super-lig/kayserispor-fenerbahce/4ztI8F3f/"
